$d = $word.ActiveDocument

# 1. Remove the duplicate visible "R.C. 2943.031. " text that precedes the
#    vanished "R.C. 2943.031" run.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "deportation, exclusion from admission into the United States, or denial of naturalization under United States law. R.C. 2943.031. "
$find.Replacement.Text = "deportation, exclusion from admission into the United States, or denial of naturalization under United States law. "
$find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)

# 2. Add a trailing space after "...following sentence:"
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Text = "accepted the plea and entered the following sentence:"
$find2.Replacement.Text = "accepted the plea and entered the following sentence: "
$find2.Execute($find2.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)
